# Build notes.docx bug-fix pass.
# Applies the following changes (see commit message: "Corrected several
# bugs on product/cart/checkout pages"):
#   1. Customer page bullet: "Sidebar needs to be identical to the company
#      page." -> "WIshlist product list needs to match design."
#   2. Remove the three Product page bullets:
#        "Add to wishlist button not working."
#        "Product reviews button not working."
#        "Write a review not working."
#   3. Move the "_GoBack" bookmark from the elastic-search bullet to the
#      start of the "Cart:" heading.
#   4. Cart bullet: "Needs a add to wishlist ..." -> "Needs an  add to
#      wishlist ..."
#   5. Remove the last two Checkout bullets:
#        "Build pop-up modal for international customers message."
#        "Hide message in shipping section when country is set to United
#        States."

$d = $word.ActiveDocument

# --- 1. Customer page bullet -------------------------------------------------
$d.Content.Find.Execute(
    "Sidebar needs to be identical to the company page.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "WIshlist product list needs to match design.", 2) | Out-Null

# --- 2. Move the _GoBack bookmark off the elastic-search bullet ------------
# It currently sits right after the close-paren run, before the trailing
# ". " run; relocate it to the very start of the "Cart:" heading paragraph
# further down (added back in step 3 below).

# --- 3. Delete the three "not working" bullets ------------------------------
# (Add to wishlist / Product reviews / Write a review) immediately after
# the elastic-search bullet and right before "Cart:".
$r1 = $d.Content
$r1.Find.Execute("Add to wishlist button not working.") | Out-Null
$delStart = $r1.Start

$r2 = $d.Content
$r2.Find.Execute("Write a review not working.") | Out-Null
$delEnd = $r2.End

$delRange = $d.Range($delStart, $delEnd)
# extend to swallow the paragraph mark of the last deleted paragraph
$delRange.MoveEnd(1, 1) | Out-Null
$delRange.Delete() | Out-Null

# --- 4. Re-add the _GoBack bookmark at the start of "Cart:" -----------------
$cartRange = $d.Content
$cartRange.Find.Execute("Cart:") | Out-Null
$cartStart = $cartRange.Start
$bmPoint = $d.Range($cartStart, $cartStart)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# --- 5. Cart bullet wording ---------------------------------------------
$d.Content.Find.Execute(
    "Needs a add to wishlist function assigned to the button AND to have it display the message before adding that item to wishlist. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Needs an  add to wishlist function assigned to the button AND to have it display the message before adding that item to wishlist. ", 2) | Out-Null

# --- 6. Delete the last two Checkout bullets --------------------------------
$r3 = $d.Content
$r3.Find.Execute("Build pop-up modal for international customers message.") | Out-Null
$delStart2 = $r3.Start

$r4 = $d.Content
$r4.Find.Execute("country is set to United States.") | Out-Null
$delEnd2 = $r4.End

$delRange2 = $d.Range($delStart2, $delEnd2)
$delRange2.MoveEnd(1, 1) | Out-Null
$delRange2.Delete() | Out-Null
